# Loan RBI, Variable Instalments
# Insert a new (blank) column N on the "Repayment schedule" sheet, shifting
# the existing Late / heading(Original) / Outstanding columns one place to
# the right (N->O, O->P, P->Q), and make "Repayment schedule" the active
# sheet/tab with a new selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this the active sheet (moves tabSelected from "Transactions" to
# "Repayment schedule" and updates the workbook's activeTab automatically).
$ws.Activate()

# Insert a blank column before column N, pushing N:P -> O:Q.
$ws.Columns("N:N").Insert()

# The newly inserted column picks up the width of its left neighbour
# (column M) when a column is inserted in Excel.
$ws.Columns("N:N").ColumnWidth = 9.8

# Update the selection on the sheet to match the new active cell.
$null = $ws.Range("K16").Select()
